$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new header cells P1, Q1 (continuing the 0..13 sequence to 14,15) ---
# Apply the same formatting used by the rest of row 1 (bold, centered, bordered header style)
$ws.Range("P1:Q1").Font.Bold = $true
$ws.Range("P1:Q1").HorizontalAlignment = -4108  # xlCenter
$ws.Range("P1:Q1").VerticalAlignment = -4160    # xlTop
$ws.Range("P1:Q1").Borders.LineStyle = 1        # xlContinuous

$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

# --- Fix existing I/K/M/O columns for rows 2-25 and add new P/Q columns ---
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 9).Value  = 2   # I -> 2
    $ws.Cells.Item($r, 11).Value = 1   # K -> 1
    $ws.Cells.Item($r, 13).Value = 2   # M -> 2
    $ws.Cells.Item($r, 15).Value = 1   # O -> 1
    $ws.Cells.Item($r, 16).Value = 2   # P -> 2 (new column)
    $ws.Cells.Item($r, 17).Value = 2   # Q -> 2 (new column)
}
